$wb = $excel.ActiveWorkbook

# ---- Delete stale work/opening rows (CmsWorkOpening, CmsWork) ----
$wsOpening = $wb.Worksheets.Item("CmsWorkOpening")
$wsOpening.Range("A2:D9").EntireRow.Delete()

$wsWork = $wb.Worksheets.Item("CmsWork")
$wsWork.Range("A2:U5").EntireRow.Delete()

# ---- SchemaOrganization ----
$ws = $wb.Worksheets.Item("SchemaOrganization")
$ws.Range("B2").Value = "http://example.com/organization1:Image0:Thumbnail200x200"
$ws.Range("B3").Value = "http://example.com/organization3:Image1:Thumbnail800x800"

# ---- CreativeCommonsLicense ----
$ws = $wb.Worksheets.Item("CreativeCommonsLicense")
$ws.Range("A3").Value = "http://creativecommons.org/licenses/by-sa/2.0/"
$ws.Range("A4").Value = "http://creativecommons.org/licenses/nc/1.0/"

# ---- CmsPerson ----
$ws = $wb.Worksheets.Item("CmsPerson")
$ws.Range("D2").Value = "http://example.com/person0:Image0"
$ws.Range("D3").Value = "http://example.com/person2:Image0:Thumbnail800x800"
$ws.Range("D4").Value = "http://example.com/person4:Image1:Thumbnail600x600"

# ---- SchemaPerson ----
$ws = $wb.Worksheets.Item("SchemaPerson")
$ws.Range("D2").Value = "http://example.com/person1:Image0:Thumbnail400x400"
$ws.Range("D3").Value = "http://example.com/person3:Image1:Thumbnail600x600"
$ws.Range("D4").Value = "http://example.com/person5:Image1:Thumbnail400x400"

# ---- CmsImage ----
$ws = $wb.Worksheets.Item("CmsImage")
$ws.Range("G6").Value = "urn:paradicms_etl:pipeline:synthetic_data:property_group:Image0:Thumbnail800x800"
$ws.Range("G24").Value = "dcterms:extent:Image0:Thumbnail400x400"
$ws.Range("G33").Value = "dcterms:language:Image0:Thumbnail200x200"
$ws.Range("G42").Value = "dcterms:medium:Image0:Thumbnail200x200"
$ws.Range("G69").Value = "dcterms:spatial:Image0:Thumbnail400x400"
$ws.Range("G78").Value = "dcterms:subject:Image0:Thumbnail200x200"
$ws.Range("G87").Value = "dcterms:title:Image0:Thumbnail200x200"
$ws.Range("G96").Value = "dcterms:type:Image0:Thumbnail800x800"
$ws.Range("G105").Value = "schema:spatial:Image0:Thumbnail400x400"
$ws.Range("G114").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:0:Image0:Thumbnail600x600"
$ws.Range("G123").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:1:Image0:Thumbnail400x400"
$ws.Range("G132").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:2:Image0:Thumbnail800x800"
$ws.Range("G141").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:3:Image0:Thumbnail400x400"
$ws.Range("G150").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:4:Image0:Thumbnail600x600"
$ws.Range("G159").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:5:Image0:Thumbnail800x800"
$ws.Range("G168").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:6:Image0:Thumbnail400x400"
$ws.Range("G195").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:9:Image0:Thumbnail600x600"
$ws.Range("G204").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:10:Image0:Thumbnail600x600"
$ws.Range("G213").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:11:Image0:Thumbnail600x600"
$ws.Range("G222").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:12:Image0:Thumbnail600x600"
$ws.Range("G240").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:14:Image0:Thumbnail800x800"
$ws.Range("G258").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:16:Image0:Thumbnail800x800"
$ws.Range("G267").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:17:Image0:Thumbnail200x200"
$ws.Range("G276").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:18:Image0:Thumbnail600x600"
$ws.Range("G285").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:19:Image0:Thumbnail600x600"
$ws.Range("G294").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:20:Image0:Thumbnail600x600"
$ws.Range("G303").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:21:Image0:Thumbnail600x600"
$ws.Range("G312").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:22:Image0:Thumbnail800x800"
$ws.Range("G321").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:23:Image0:Thumbnail200x200"
$ws.Range("G330").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:24:Image0:Thumbnail600x600"
$ws.Range("G348").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:26:Image0:Thumbnail400x400"
$ws.Range("G366").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:28:Image0:Thumbnail200x200"
$ws.Range("G384").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:30:Image0:Thumbnail800x800"
$ws.Range("G393").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:31:Image0:Thumbnail200x200"
$ws.Range("G402").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:32:Image0:Thumbnail800x800"
$ws.Range("G411").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:33:Image0:Thumbnail200x200"
$ws.Range("G420").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:34:Image0:Thumbnail600x600"
$ws.Range("G429").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:35:Image0:Thumbnail600x600"
$ws.Range("G447").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:37:Image0:Thumbnail200x200"
$ws.Range("G465").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:39:Image0:Thumbnail400x400"
$ws.Range("G474").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:40:Image0:Thumbnail800x800"
$ws.Range("G483").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:41:Image0:Thumbnail800x800"
$ws.Range("G519").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:45:Image0:Thumbnail400x400"
$ws.Range("G528").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:46:Image0:Thumbnail200x200"
$ws.Range("G537").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:47:Image0:Thumbnail400x400"
$ws.Range("G564").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:50:Image0:Thumbnail400x400"
$ws.Range("G591").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:53:Image0:Thumbnail200x200"
$ws.Range("G600").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:54:Image0:Thumbnail600x600"
$ws.Range("G609").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:55:Image0:Thumbnail400x400"
$ws.Range("G627").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:57:Image0:Thumbnail200x200"
$ws.Range("G645").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:59:Image0:Thumbnail400x400"
$ws.Range("G654").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:60:Image0:Thumbnail600x600"
$ws.Range("G663").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:61:Image0:Thumbnail600x600"
$ws.Range("G672").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:62:Image0:Thumbnail800x800"
$ws.Range("G681").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:63:Image0:Thumbnail800x800"
$ws.Range("G690").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:64:Image0:Thumbnail800x800"
$ws.Range("G699").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:65:Image0:Thumbnail200x200"
$ws.Range("G708").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:66:Image0:Thumbnail200x200"
$ws.Range("G735").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:69:Image0:Thumbnail400x400"
$ws.Range("G753").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:71:Image0:Thumbnail600x600"
$ws.Range("G771").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:73:Image0:Thumbnail600x600"
$ws.Range("G780").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:74:Image0:Thumbnail200x200"
$ws.Range("G789").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:75:Image0:Thumbnail400x400"
$ws.Range("G798").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:76:Image0:Thumbnail400x400"
$ws.Range("G807").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:77:Image0:Thumbnail600x600"
$ws.Range("G816").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:78:Image0:Thumbnail800x800"
$ws.Range("G825").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:79:Image0:Thumbnail800x800"
$ws.Range("G834").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:80:Image0:Thumbnail800x800"
$ws.Range("G843").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:81:Image0:Thumbnail600x600"
$ws.Range("G852").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:82:Image0:Thumbnail800x800"
$ws.Range("G861").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:83:Image0:Thumbnail200x200"
$ws.Range("G870").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:84:Image0:Thumbnail800x800"
$ws.Range("G888").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:86:Image0:Thumbnail400x400"
$ws.Range("G897").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:87:Image0:Thumbnail800x800"
$ws.Range("G906").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:88:Image0:Thumbnail200x200"
$ws.Range("G915").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:89:Image0:Thumbnail800x800"
$ws.Range("G924").Value = "http://example.com/organization0:Image0:Thumbnail600x600"
$ws.Range("G942").Value = "http://example.com/organization2:Image0:Thumbnail800x800"
$ws.Range("G951").Value = "http://example.com/organization3:Image0:Thumbnail200x200"
$ws.Range("G960").Value = "http://example.com/organization4:Image0:Thumbnail200x200"
$ws.Range("G969").Value = "http://example.com/organization5:Image0:Thumbnail600x600"
$ws.Range("G978").Value = "http://example.com/person0:Image0:Thumbnail400x400"
$ws.Range("G988").Value = "http://example.com/person1:Image0:Thumbnail200x200"
$ws.Range("G997").Value = "http://example.com/person2:Image0:Thumbnail600x600"
$ws.Range("G1006").Value = "http://example.com/person3:Image0:Thumbnail800x800"
$ws.Range("G1015").Value = "http://example.com/person4:Image0:Thumbnail600x600"
$ws.Range("G1024").Value = "http://example.com/person5:Image0:Thumbnail600x600"
$ws.Range("G1033").Value = "http://example.com/collection0/work0:Image0:Thumbnail800x800"
$ws.Range("G1043").Value = "http://example.com/collection0/work1:Image0:Thumbnail400x400"
$ws.Range("G1052").Value = "http://example.com/collection0/work2:Image0:Thumbnail800x800"
$ws.Range("G1061").Value = "http://example.com/collection0/work3:Image0:Thumbnail600x600"
$ws.Range("G1070").Value = "http://example.com/collection1:Image0:Thumbnail200x200"
$ws.Range("G1079").Value = "http://example.com/collection1/work4:Image0:Thumbnail400x400"
$ws.Range("G1106").Value = "http://example.com/collection1/work7:Image0:Thumbnail800x800"
$ws.Range("G1124").Value = "http://example.com/freestandingwork9:Image0:Thumbnail600x600"
$ws.Range("G1133").Value = "http://example.com/freestandingwork10:Image0:Thumbnail600x600"
$ws.Range("G1142").Value = "http://example.com/freestandingwork11:Image0:Thumbnail400x400"

# ---- CmsProperty ----
$ws = $wb.Worksheets.Item("CmsProperty")
$ws.Range("C2").Value = "dcterms:description:Image0"
$ws.Range("C3").Value = "dcterms:extent:Image0:Thumbnail400x400"
$ws.Range("C4").Value = "dcterms:language:Image0:Thumbnail800x800"
$ws.Range("C5").Value = "dcterms:medium:Image1:Thumbnail800x800"
$ws.Range("C6").Value = "dcterms:publisher:Image1:Thumbnail800x800"
$ws.Range("C7").Value = "dcterms:source:Image1:Thumbnail400x400"
$ws.Range("C8").Value = "dcterms:spatial:Image0:Thumbnail200x200"
$ws.Range("C9").Value = "dcterms:subject:Image0:Thumbnail200x200"
$ws.Range("C10").Value = "dcterms:title:Image1"
$ws.Range("C11").Value = "dcterms:type:Image1:Thumbnail600x600"

# ---- SchemaProperty ----
$ws = $wb.Worksheets.Item("SchemaProperty")
$ws.Range("C2").Value = "https://schema.org/spatial:Image1:Thumbnail600x600"

# ---- CmsPropertyGroup ----
$ws = $wb.Worksheets.Item("CmsPropertyGroup")
$ws.Range("C2").Value = "urn:paradicms_etl:pipeline:synthetic_data:property_group:Image1:Thumbnail600x600"

# ---- CmsConcept ----
$ws = $wb.Worksheets.Item("CmsConcept")
$ws.Range("B2").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:0:Image0"
$ws.Range("B4").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:2:Image0:Thumbnail400x400"
$ws.Range("B5").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:3:Image0:Thumbnail200x200"
$ws.Range("B6").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:4:Image0"
$ws.Range("B7").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:5:Image1:Thumbnail800x800"
$ws.Range("B8").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:6:Image1"
$ws.Range("B9").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:7:Image0:Thumbnail600x600"
$ws.Range("B10").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:8:Image1:Thumbnail800x800"
$ws.Range("B12").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:10:Image1:Thumbnail200x200"
$ws.Range("B13").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:11:Image1"
$ws.Range("B14").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:12:Image0:Thumbnail400x400"
$ws.Range("B16").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:14:Image0:Thumbnail800x800"
$ws.Range("B17").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:15:Image1:Thumbnail800x800"
$ws.Range("B19").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:17:Image1:Thumbnail200x200"
$ws.Range("B20").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:18:Image1:Thumbnail400x400"
$ws.Range("B21").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:19:Image1:Thumbnail200x200"
$ws.Range("B22").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:20:Image1:Thumbnail200x200"
$ws.Range("B23").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:21:Image0"
$ws.Range("B25").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:23:Image1:Thumbnail800x800"
$ws.Range("B26").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:24:Image0:Thumbnail800x800"
$ws.Range("B27").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:25:Image1:Thumbnail600x600"
$ws.Range("B28").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:26:Image0:Thumbnail400x400"
$ws.Range("B30").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:28:Image0:Thumbnail400x400"
$ws.Range("B31").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:29:Image1:Thumbnail600x600"
$ws.Range("B32").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:30:Image1:Thumbnail800x800"
$ws.Range("B33").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:31:Image0:Thumbnail800x800"
$ws.Range("B34").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:32:Image0:Thumbnail800x800"
$ws.Range("B35").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:33:Image1:Thumbnail200x200"
$ws.Range("B36").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:34:Image0:Thumbnail200x200"
$ws.Range("B37").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:35:Image0:Thumbnail600x600"
$ws.Range("B38").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:36:Image1:Thumbnail200x200"
$ws.Range("B39").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:37:Image1:Thumbnail400x400"
$ws.Range("B41").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:39:Image1:Thumbnail600x600"
$ws.Range("B42").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:40:Image0"
$ws.Range("B43").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:41:Image0:Thumbnail600x600"
$ws.Range("B44").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:42:Image0:Thumbnail400x400"
$ws.Range("B45").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:43:Image0"
$ws.Range("B46").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:44:Image1:Thumbnail600x600"
$ws.Range("B47").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:45:Image0:Thumbnail200x200"
$ws.Range("B48").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:46:Image1:Thumbnail800x800"
$ws.Range("B49").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:47:Image0"
$ws.Range("B50").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:48:Image0:Thumbnail400x400"
$ws.Range("B51").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:49:Image0"
$ws.Range("B52").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:50:Image1:Thumbnail600x600"
$ws.Range("B53").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:51:Image0:Thumbnail200x200"
$ws.Range("B54").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:52:Image0:Thumbnail600x600"
$ws.Range("B55").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:53:Image1:Thumbnail800x800"
$ws.Range("B56").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:54:Image1:Thumbnail400x400"
$ws.Range("B57").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:55:Image1:Thumbnail400x400"
$ws.Range("B58").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:56:Image0:Thumbnail400x400"
$ws.Range("B59").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:57:Image1:Thumbnail400x400"
$ws.Range("B60").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:58:Image1"
$ws.Range("B61").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:59:Image1:Thumbnail200x200"
$ws.Range("B62").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:60:Image1:Thumbnail400x400"
$ws.Range("B64").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:62:Image0"
$ws.Range("B65").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:63:Image1:Thumbnail200x200"
$ws.Range("B66").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:64:Image1"
$ws.Range("B67").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:65:Image0:Thumbnail200x200"
$ws.Range("B68").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:66:Image0:Thumbnail400x400"
$ws.Range("B69").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:67:Image0:Thumbnail200x200"
$ws.Range("B70").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:68:Image0:Thumbnail600x600"
$ws.Range("B71").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:69:Image1:Thumbnail400x400"
$ws.Range("B73").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:71:Image0"
$ws.Range("B74").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:72:Image0:Thumbnail400x400"
$ws.Range("B75").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:73:Image0:Thumbnail400x400"
$ws.Range("B76").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:74:Image1:Thumbnail800x800"
$ws.Range("B77").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:75:Image1:Thumbnail200x200"
$ws.Range("B78").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:76:Image1"
$ws.Range("B79").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:77:Image0:Thumbnail400x400"
$ws.Range("B80").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:78:Image1:Thumbnail800x800"
$ws.Range("B81").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:79:Image0:Thumbnail400x400"

# ---- SchemaDefinedTerm ----
$ws = $wb.Worksheets.Item("SchemaDefinedTerm")
$ws.Range("B3").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:81:Image1:Thumbnail400x400"
$ws.Range("B4").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:82:Image0:Thumbnail600x600"
$ws.Range("B5").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:83:Image1:Thumbnail200x200"
$ws.Range("B6").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:84:Image0:Thumbnail600x600"
$ws.Range("B7").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:85:Image0:Thumbnail800x800"
$ws.Range("B8").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:86:Image1:Thumbnail800x800"
$ws.Range("B9").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:87:Image1"
$ws.Range("B11").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:89:Image0:Thumbnail800x800"

# ---- CmsOrganization ----
$ws = $wb.Worksheets.Item("CmsOrganization")
$ws.Range("B2").Value = "http://example.com/organization0:Image0:Thumbnail800x800"
$ws.Range("B3").Value = "http://example.com/organization2:Image1:Thumbnail200x200"
$ws.Range("B4").Value = "http://example.com/organization4:Image0"
